$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.552.18"
$ws.Range("D3").Value = "1.751.59"
$ws.Range("E3").Value = "  -0.27%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4492"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -1.97%  "
$ws.Range("E9").Value = "  -1.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("E11").Value = "  -2.18%  "
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.975"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.142"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.23%  "
$ws.Range("D16").Value = "1.752.56"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.98%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001054"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06379"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.78%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.726"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.40%  "
$ws.Range("D23").Value = "27.600.71"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.084"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.37%  "
$ws.Range("E27").Value = "  -1.83%  "
$ws.Range("D28").Value = "1.953.15"
$ws.Range("E28").Value = "  -0.29%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.099"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.68%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "125.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.081"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09157"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.481"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.28%  "
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2087"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05996"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.6263"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.906"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("E41").Value = "  -1.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.395"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.766"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.714"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5841"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.27%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "121.71"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.928"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.06878"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.124"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.58%  "
